$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 13.5642721015424
$ws.Range("C2").Value = 8.773884001667684
$ws.Range("D2").Value = 10.75636708064979
$ws.Range("F2").Value = 30.36434137114557
$ws.Range("G2").Value = 3.633931526186452
$ws.Range("I2").Value = 20.03760700138333
$ws.Range("J2").Value = 11.07709599258071
$ws.Range("M2").Value = 17.55538640893379
$ws.Range("O2").Value = 22.02576866545462
$ws.Range("B3").Value = 12.93852324362643
$ws.Range("C3").Value = 8.249131485747323
$ws.Range("D3").Value = 10.73558238273878
$ws.Range("F3").Value = 30.45757176190199
$ws.Range("G3").Value = 3.636037852679719
$ws.Range("I3").Value = 20.18423156800726
$ws.Range("J3").Value = 11.11253272763376
$ws.Range("M3").Value = 17.34338969366835
$ws.Range("O3").Value = 22.12841502399705
$ws.Range("B4").Value = 12.53863524316537
$ws.Range("C4").Value = 7.9087958618406
$ws.Range("D4").Value = 10.72437346580099
$ws.Range("F4").Value = 30.52413094511854
$ws.Range("G4").Value = 3.637399453174576
$ws.Range("I4").Value = 20.2794502834624
$ws.Range("J4").Value = 11.13598152408213
$ws.Range("M4").Value = 17.21404259963569
$ws.Range("O4").Value = 22.19780755316026
$ws.Range("B5").Value = 12.37192891948538
$ws.Range("C5").Value = 7.765607481807205
$ws.Range("D5").Value = 10.72019935986136
$ws.Range("F5").Value = 30.55358679790189
$ws.Range("G5").Value = 3.637971547835957
$ws.Range("I5").Value = 20.31955811660882
$ws.Range("J5").Value = 11.14596219508651
$ws.Range("M5").Value = 17.16159095627162
$ws.Range("O5").Value = 22.22768026138867
$ws.Range("B6").Value = 12.34402727286292
$ws.Range("C6").Value = 7.741561549195625
$ws.Range("D6").Value = 10.71953011204599
$ws.Range("F6").Value = 30.55861848951791
$ws.Range("G6").Value = 3.638067586068606
$ws.Range("I6").Value = 20.32629684460854
$ws.Range("J6").Value = 11.14764515294494
$ws.Range("M6").Value = 17.15289852362582
$ws.Range("O6").Value = 22.23273672796764
$ws.Range("B7").Value = 12.53640188503537
$ws.Range("C7").Value = 7.906882896290357
$ws.Range("D7").Value = 10.72431557468555
$ws.Range("F7").Value = 30.52451876676746
$ws.Range("G7").Value = 3.637407098801518
$ws.Range("I7").Value = 20.27998590575954
$ws.Range("J7").Value = 11.1361144056958
$ws.Range("M7").Value = 17.21333410366548
$ws.Range("O7").Value = 22.19820397868369
$ws.Range("B8").Value = 13.35187432693163
$ws.Range("C8").Value = 8.596745601145997
$ws.Range("D8").Value = 10.74888014020064
$ws.Range("F8").Value = 30.39454800714731
$ws.Range("G8").Value = 3.63464364346159
$ws.Range("I8").Value = 20.08708594303341
$ws.Range("J8").Value = 11.08896372048423
$ws.Range("M8").Value = 17.48215120509374
$ws.Range("O8").Value = 22.05983578528188
$ws.Range("B9").Value = 14.8195494109944
$ws.Range("C9").Value = 9.803672901091991
$ws.Range("D9").Value = 10.80922646344669
$ws.Range("F9").Value = 30.21400574971397
$ws.Range("G9").Value = 3.629763996386234
$ws.Range("I9").Value = 19.74998624544916
$ws.Range("J9").Value = 11.00991440187285
$ws.Range("M9").Value = 18.01350612834174
$ws.Range("O9").Value = 21.83929955188875
$ws.Range("B10").Value = 15.80992006425195
$ws.Range("C10").Value = 10.59957161222633
$ws.Range("D10").Value = 10.86077389916139
$ws.Range("F10").Value = 30.12720042557575
$ws.Range("G10").Value = 3.626504285422134
$ws.Range("I10").Value = 19.52741031123113
$ws.Range("J10").Value = 10.96001336041442
$ws.Range("M10").Value = 18.40339226489687
$ws.Range("O10").Value = 21.70862166814443
$ws.Range("B11").Value = 16.24014158462237
$ws.Range("C11").Value = 10.94169480526663
$ws.Range("D11").Value = 10.88573967586871
$ws.Range("F11").Value = 30.09776025235734
$ws.Range("G11").Value = 3.625091260104659
$ws.Range("I11").Value = 19.43160307799745
$ws.Range("J11").Value = 10.93908756037013
$ws.Range("M11").Value = 18.58001389680982
$ws.Range("O11").Value = 21.65605894445671
$ws.Range("B12").Value = 16.40005540702944
$ws.Range("C12").Value = 11.06836971328369
$ws.Range("D12").Value = 10.89540696260404
$ws.Range("F12").Value = 30.08806342483125
$ws.Range("G12").Value = 3.624566169430219
$ws.Range("I12").Value = 19.39610663122658
$ws.Range("J12").Value = 10.93141871125418
$ws.Range("M12").Value = 18.64673736048536
$ws.Range("O12").Value = 21.63715090203791
$ws.Range("B13").Value = 16.36574976497187
$ws.Range("C13").Value = 11.04121623892484
$ws.Range("D13").Value = 10.89331553222545
$ws.Range("F13").Value = 30.09008715722421
$ws.Range("G13").Value = 3.624678813536139
$ws.Range("I13").Value = 19.40371655656875
$ws.Range("J13").Value = 10.93305897852114
$ws.Range("M13").Value = 18.6323751834104
$ws.Range("O13").Value = 21.6411786738926
$ws.Range("B14").Value = 16.25335837036196
$ws.Range("C14").Value = 10.95217419783783
$ws.Range("D14").Value = 10.88653076304704
$ws.Range("F14").Value = 30.09693336302067
$ws.Range("G14").Value = 3.625047860626211
$ws.Range("I14").Value = 19.4286670524449
$ws.Range("J14").Value = 10.93845152221432
$ws.Range("M14").Value = 18.58550676787724
$ws.Range("O14").Value = 21.6544833640265
$ws.Range("B15").Value = 16.18412228709618
$ws.Range("C15").Value = 10.89725804689561
$ws.Range("D15").Value = 10.88240252616598
$ws.Range("F15").Value = 30.10131607467914
$ws.Range("G15").Value = 3.625275212296587
$ws.Range("I15").Value = 19.44405203747239
$ws.Range("J15").Value = 10.9417878622423
$ws.Range("M15").Value = 18.5567762131725
$ws.Range("O15").Value = 21.66276281140691
$ws.Range("B16").Value = 15.7813881699721
$ws.Range("C16").Value = 10.57681068012397
$ws.Range("D16").Value = 10.85917244020725
$ws.Range("F16").Value = 30.12932718986011
$ws.Range("G16").Value = 3.626598030842554
$ws.Range("I16").Value = 19.53378115187562
$ws.Range("J16").Value = 10.96141662834899
$ws.Range("M16").Value = 18.39183018673513
$ws.Range("O16").Value = 21.71219583582601
$ws.Range("B17").Value = 15.52906057724642
$ws.Range("C17").Value = 10.3751098211428
$ws.Range("D17").Value = 10.84530663501122
$ws.Range("F17").Value = 30.14908964777544
$ws.Range("G17").Value = 3.627427387074775
$ws.Range("I17").Value = 19.59022189461407
$ws.Range("J17").Value = 10.97391281138553
$ws.Range("M17").Value = 18.29041363619292
$ws.Range("O17").Value = 21.74428906021198
$ws.Range("B18").Value = 15.38202130953637
$ws.Range("C18").Value = 10.25722086034779
$ws.Range("D18").Value = 10.83747444169824
$ws.Range("F18").Value = 30.16140191641303
$ws.Range("G18").Value = 3.627910986910455
$ws.Range("I18").Value = 19.62319744024137
$ws.Range("J18").Value = 10.98126729515844
$ws.Range("M18").Value = 18.23201500470381
$ws.Range("O18").Value = 21.76339574397645
$ws.Range("B19").Value = 15.33191153590602
$ws.Range("C19").Value = 10.21698396634987
$ws.Range("D19").Value = 10.8348473069316
$ws.Range("F19").Value = 30.16573279965493
$ws.Range("G19").Value = 3.628075856501072
$ws.Range("I19").Value = 19.63445036093205
$ws.Range("J19").Value = 10.98378607700876
$ws.Range("M19").Value = 18.21223244285963
$ws.Range("O19").Value = 21.76997594811281
$ws.Range("B20").Value = 15.55611929998258
$ws.Range("C20").Value = 10.39677548710338
$ws.Range("D20").Value = 10.84676790236894
$ws.Range("F20").Value = 30.14688800138182
$ws.Range("G20").Value = 3.627338420457427
$ws.Range("I20").Value = 19.58416065171318
$ws.Range("J20").Value = 10.97256528550168
$ws.Range("M20").Value = 18.30121687295227
$ws.Range("O20").Value = 21.74080561910083
$ws.Range("B21").Value = 16.28645249789735
$ws.Range("C21").Value = 10.97840623879457
$ws.Range("D21").Value = 10.88851786501077
$ws.Range("F21").Value = 30.09488302407459
$ws.Range("G21").Value = 3.62493919178927
$ws.Range("I21").Value = 19.42131721062502
$ws.Range("J21").Value = 10.93686067162554
$ws.Range("M21").Value = 18.5992778968906
$ws.Range("O21").Value = 21.65054836317145
$ws.Range("B22").Value = 16.74624258559598
$ws.Range("C22").Value = 11.34174749261378
$ws.Range("D22").Value = 10.91704470076279
$ws.Range("F22").Value = 30.06935779670732
$ws.Range("G22").Value = 3.623429371143551
$ws.Range("I22").Value = 19.31945796422199
$ws.Range("J22").Value = 10.91501377486398
$ws.Range("M22").Value = 18.79312581527289
$ws.Range("O22").Value = 21.59737023799696
$ws.Range("B23").Value = 16.50247105306647
$ws.Range("C23").Value = 11.14936465099731
$ws.Range("D23").Value = 10.90170750124645
$ws.Range("F23").Value = 30.08220477036508
$ws.Range("G23").Value = 3.624229880885846
$ws.Range("I23").Value = 19.37340382890259
$ws.Range("J23").Value = 10.9265376650545
$ws.Range("M23").Value = 18.68976980205379
$ws.Range("O23").Value = 21.62521866539945
$ws.Range("B24").Value = 15.54389218981617
$ws.Range("C24").Value = 10.38698644522684
$ws.Range("D24").Value = 10.84610682869315
$ws.Range("F24").Value = 30.14788040523541
$ws.Range("G24").Value = 3.627378621104789
$ws.Range("I24").Value = 19.5868992980392
$ws.Range("J24").Value = 10.97317397160445
$ws.Range("M24").Value = 18.29633301494006
$ws.Range("O24").Value = 21.74237844031412
$ws.Range("B25").Value = 14.43746726912631
$ws.Range("C25").Value = 9.493016839034226
$ws.Range("D25").Value = 10.79161851088895
$ws.Range("F25").Value = 30.25483236162476
$ws.Range("G25").Value = 3.631026679963398
$ws.Range("I25").Value = 19.83677402807435
$ws.Range("J25").Value = 11.02986344218237
$ws.Range("M25").Value = 17.86962845686895
$ws.Range("O25").Value = 21.8934819716218
